# PlayerData.xlsx update
# - Adds new "Plus" bonus stat columns (MaxHPPlus, MaxMPPlus, CarryWeightPlus,
#   MaxStaminaPlus, MaxSatietyPlus, StrengthPlus, EndurancePlus, DexterityPlus,
#   PerceptionPlus, MasterPlu, WillPlus, Magical_powerplus, CharismaPlus,
#   AgilityPlus, LuckPlus) interleaved with the existing stat columns and
#   reorders EXP/NextLevelEXP/Satiety/Perception, extending the sheet from
#   A:Y out to A:AN.
# - Refreshes the row 3 "Player" template values to include starting bonus
#   stats (damage_create system fix).
# - Refreshes the selected cell in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 1), column A..AN
$headerNames = @(
    "ID", "Name", "Level", "MaxHP", "MaxHPPlus", "HP", "MaxMP", "MaxMPPlus", "MP",
    "CarryWeight", "CarryWeightPlus", "Weight", "MaxStamina", "MaxStaminaPlus", "Stamina",
    "MaxSatiety", "MaxSatietyPlus", "Satiety", "EXP", "NextLevelEXP",
    "Strength", "StrengthPlus", "Endurance", "EndurancePlus", "Dexterity", "DexterityPlus",
    "Perception", "PerceptionPlus", "Master", "MasterPlu", "Will", "WillPlus",
    "Magical_power", "Magical_powerplus", "Charisma", "CharismaPlus",
    "Agility", "AgilityPlus", "Luck", "LuckPlus"
)

# Row 2 ("None" placeholder row) values, column A..AN
$row2Values = @(
    0, "None", 0, 0, 0, 0, 0, 0, 0, 0,
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0,
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0,
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0
)

# Row 3 ("Player" template) values, column A..AN
$row3Values = @(
    1, "Player", 1, 100, 0, 100, 100, 0, 100, 100,
    0, 0, 100, 0, 0, 0, 0, 0, 50, 5,
    0, 5, 0, 5, 0, 3, 0, 3, 0, 3,
    0, 1, 0, 2, 0, 1, 0, 3, 0, 0
)

$colCount = $headerNames.Length

for ($i = 0; $i -lt $colCount; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headerNames[$i]
    $ws.Cells.Item(2, $col).Value = $row2Values[$i]
    $ws.Cells.Item(3, $col).Value = $row3Values[$i]
}

# Update the sheet view selection/scroll state to match the saved workbook
$ws.Range("K14").Select()
